$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking prices (e.g. "218.10") keep
# their literal formatting instead of being parsed into Doubles.
$ws.Range("D2:D51").NumberFormat = "@"

# Price (D) updates
$ws.Range("D2").Value = "26.253.99"
$ws.Range("D3").Value = "1.674.30"
$ws.Range("D5").Value = "218.10"
$ws.Range("D6").Value = "0.5135"
$ws.Range("D7").Value = "1.006"
$ws.Range("D8").Value = "0.2663"
$ws.Range("D9").Value = "0.06409"
$ws.Range("D10").Value = "21.57"
$ws.Range("D11").Value = "0.07374"
$ws.Range("D12").Value = "1.672.37"
$ws.Range("D13").Value = "4.556"
$ws.Range("D14").Value = "0.5829"
$ws.Range("D15").Value = "1.900.74"
$ws.Range("D16").Value = "0.000008698"
$ws.Range("D17").Value = "64.91"
$ws.Range("D18").Value = "26.329.20"
$ws.Range("D19").Value = "4.961"
$ws.Range("D21").Value = "10.85"
$ws.Range("D22").Value = "190.10"
$ws.Range("D23").Value = "6.226"
$ws.Range("D24").Value = "1.007"
$ws.Range("D25").Value = "144.33"
$ws.Range("D26").Value = "7.632"
$ws.Range("D27").Value = "0.1184"
$ws.Range("D28").Value = "15.64"
$ws.Range("D29").Value = "0.05934"
$ws.Range("D30").Value = "1.282"
$ws.Range("D32").Value = "3.536"
$ws.Range("D33").Value = "3.526"
$ws.Range("D35").Value = "1.015"
$ws.Range("D36").Value = "0.6020"
$ws.Range("D37").Value = "2.367"
$ws.Range("D38").Value = "2.652"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D40").Value = "6.065"
$ws.Range("D41").Value = "1.078.07"
$ws.Range("D42").Value = "0.8693"
$ws.Range("D44").Value = "99.90"
$ws.Range("D45").Value = "1.822.11"
$ws.Range("D46").Value = "0.00000000114"
$ws.Range("D47").Value = "56.02"
$ws.Range("D49").Value = "8.065"
$ws.Range("D50").Value = "0.4300"
$ws.Range("D51").Value = "0.05204"

# Volume (E) updates
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  -1.11%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  +5.45%  "
$ws.Range("E9").Value = "  +5.47%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +9.08%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("E29").Value = "  +3.78%  "
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("E37").Value = "  -3.10%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("E41").Value = "  +1.96%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  +4.75%  "
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("E46").Value = "  +5.76%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("E51").Value = "  -0.31%  "

# Restore the original (unstyled) cell style now that the text values are set,
# so the saved workbook has no residual style/number-format delta.
$ws.Range("D2:D51").Style = "Normal"
